# Updated symbol list (price/volume refresh) applied via Excel COM interop.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "E2" "-0.54%"
Set-TextValue "D3" "28.42"
Set-TextValue "E3" "-4.59%"
Set-TextValue "D4" "5.239"
Set-TextValue "E4" "1.28%"
Set-TextValue "D5" "0.05699"
Set-TextValue "E5" "-0.53%"
Set-TextValue "D6" "6.611"
Set-TextValue "D7" "3.195"
Set-TextValue "E7" "3.17%"
Set-TextValue "D8" "0.8499"
Set-TextValue "E8" "-0.73%"
Set-TextValue "D9" "0.8835"
Set-TextValue "E9" "1.84%"
Set-TextValue "D10" "0.1367"
Set-TextValue "E10" "0.06%"
Set-TextValue "D11" "0.07022"
Set-TextValue "E11" "-0.76%"
Set-TextValue "D12" "0.03138"
Set-TextValue "E12" "7.12%"
Set-TextValue "D13" "0.09209"
Set-TextValue "E13" "-1.84%"
Set-TextValue "D14" "0.001535"
Set-TextValue "E14" "1.51%"
Set-TextValue "D15" "0.0005966"
Set-TextValue "E15" "-0.67%"
Set-TextValue "D16" "0.005967"
Set-TextValue "E16" "-3.42%"
Set-TextValue "D17" "3.492"
Set-TextValue "E17" "0.10%"
Set-TextValue "E18" "-0.50%"
Set-TextValue "D19" "0.3169"
Set-TextValue "E19" "0.44%"
Set-TextValue "D20" "0.03283"
Set-TextValue "E20" "-4.50%"
Set-TextValue "D21" "0.1288"
Set-TextValue "E21" "-1.17%"
Set-TextValue "D22" "3.509"
Set-TextValue "E22" "1.33%"
Set-TextValue "D23" "0.04091"
Set-TextValue "E23" "-1.21%"
Set-TextValue "D24" "0.1378"
Set-TextValue "E24" "-0.08%"
Set-TextValue "D25" "0.001220"
Set-TextValue "E25" "-0.11%"
Set-TextValue "D26" "0.004145"
Set-TextValue "E26" "-17.34%"
Set-TextValue "E27" "-0.83%"
Set-TextValue "D28" "0.0001448"
Set-TextValue "D40" "0.03773"
Set-TextValue "E40" "0.52%"
Set-TextValue "E41" "-0.84%"
Set-TextValue "D42" "0.003737"
Set-TextValue "E42" "7.09%"
Set-TextValue "E43" "-10.21%"
Set-TextValue "D44" "0.009186"
Set-TextValue "E44" "-4.01%"
Set-TextValue "D45" "0.00005272"
Set-TextValue "E45" "0.27%"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "E46" "-0.01%"
Set-TextValue "D47" "0.1049"
Set-TextValue "E47" "62.27%"
Set-TextValue "D48" "0.002437"
Set-TextValue "E48" "-3.61%"
Set-TextValue "E49" "-0.01%"
Set-TextValue "E50" "-0.01%"
